$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ELEGANCE LASTING LIP STICK 6 -> ELEGANCE LINE SHADOW POWDER 2
$ws.Range("A2").Value = 19
$ws.Range("B2").Value = 2987
$ws.Range("C2").Value = "ELEGANCE LINE SHADOW POWDER 2"
$ws.Range("D2").Value = "ELEGANCE LINE SHADOW POWDER 2"
$ws.Range("I2").Value = "ELEGANCE LINE SHADOW POWDER 2"
$ws.Range("J2").Value = "ELEGANCE LINE SHADOW POWDER 2"

# Row 3: ELEGANCE LINE SHADOW POWDER 1 -> ELEGANCE LINE SHADOW POWDER 3
$ws.Range("B3").Value = 2988
$ws.Range("C3").Value = "ELEGANCE LINE SHADOW POWDER 3"
$ws.Range("D3").Value = "ELEGANCE LINE SHADOW POWDER 3"
$ws.Range("I3").Value = "ELEGANCE LINE SHADOW POWDER 3"
$ws.Range("J3").Value = "ELEGANCE LINE SHADOW POWDER 3"
